# Add the new "creditCalcCard" worksheet after the existing "makeDeposite"
# sheet, populate its header/label cells (which land in shared strings),
# and give those cells the new Calibri-based style.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "creditCalcCard"

$ws2.Range("A1").Value = "tab1"
$ws2.Range("B1").Value = "tab2"
$ws2.Range("C1").Value = "tab3"
$ws2.Range("A2").Value = "Platinum картка100&nbsp;000"
$ws2.Range("C2").Value = "9 місяців"

$ws2.Range("A1").Font.Name = "Calibri"
$ws2.Range("B1").Font.Name = "Calibri"
$ws2.Range("C1").Font.Name = "Calibri"
$ws2.Range("A2").Font.Name = "Calibri"
$ws2.Range("C2").Font.Name = "Calibri"
